# Apply the table-style change found on slide 5:
#   <a:tableStyleId>{7BFD4734-97B8-4E0E-8FBF-377300898EE9}</a:tableStyleId>
#   -> {45CB31FE-C1D9-4707-8F23-C2FF372FE7D4}
#
# (PowerPoint's object model has no "Table.Style = ..." setter - the
# correct call is Table.ApplyStyle("{GUID}"), which rewrites the
# <a:tableStyleId> element on the shape's <a:tbl>.)

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(5)

for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $shape = $s.Shapes.Item($i)
    if ($shape.HasTable) {
        $shape.Table.ApplyStyle("{45CB31FE-C1D9-4707-8F23-C2FF372FE7D4}")
    }
}
